$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customers")
$tbl = $ws.ListObjects.Item(1)

$newCol = $tbl.ListColumns.Add()

$headerCell = $ws.Range("D1")
$headerCell.NumberFormat = "@"
$headerCell.Value = "Updated Phone No"

$dataCell = $ws.Range("D2")
$dataCell.NumberFormat = "@"
$dataCell.Value = "9038655199"
